$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for Actual Consumption / Timestamp, rows 2..41 (row 1 is header)
$data = @(
    @(4622, 45810),
    @(4520, 45810.01041666666),
    @(4502, 45810.02083333334),
    @(4573, 45810.03125),
    @(4525, 45810.04166666666),
    @(4515, 45810.05208333334),
    @(4463, 45810.0625),
    @(4434, 45810.07291666666),
    @(4408, 45810.08333333334),
    @(4416, 45810.09375),
    @(4433, 45810.10416666666),
    @(4427, 45810.11458333334),
    @(4380, 45810.125),
    @(4478, 45810.13541666666),
    @(4459, 45810.14583333334),
    @(4465, 45810.15625),
    @(4524, 45810.16666666666),
    @(4501, 45810.17708333334),
    @(4535, 45810.1875),
    @(4687, 45810.19791666666),
    @(4870, 45810.20833333334),
    @(4971, 45810.21875),
    @(5059, 45810.22916666666),
    @(5215, 45810.23958333334),
    @(5456, 45810.25),
    @(5530, 45810.26041666666),
    @(5633, 45810.27083333334),
    @(5622, 45810.28125),
    @(5603, 45810.29166666666),
    @(5617, 45810.30208333334),
    @(5656, 45810.3125),
    @(5616, 45810.32291666666),
    @(5520, 45810.33333333334),
    @(5472, 45810.34375),
    @(5406, 45810.35416666666),
    @(5429, 45810.36458333334),
    @(5265, 45810.375),
    @(5098, 45810.38541666666),
    @(5087, 45810.39583333334),
    @(5040, 45810.40625),
)

$startRow = 2
$lastExistingRow = 27

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]

    # Newly created rows (28+) don't inherit the timestamp cell's number
    # format automatically - apply it explicitly so it matches the rest
    # of column B.
    if ($r -gt $lastExistingRow) {
        $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    }
}
